$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B4").Value = -0.55
$summary.Range("B5").Value = -0.21
$summary.Range("B6").Value = 53
$summary.Range("B7").Value = 19
$summary.Range("B9").Value = 35.85

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 53
$status.Range("E4").Value = -0.55
$status.Range("G4").Value = 35.85

# --- Append new closed trade row to "All Trades" and "MarketMaking" sheets ---
$newRow = 54

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item($newRow, 1).Value = 53

    $ws.Cells.Item($newRow, 2).NumberFormat = "@"
    $ws.Cells.Item($newRow, 2).Value = "2026-02-17"
    $ws.Cells.Item($newRow, 2).Style = "Normal"

    $ws.Cells.Item($newRow, 3).NumberFormat = "@"
    $ws.Cells.Item($newRow, 3).Value = "08:42:14"
    $ws.Cells.Item($newRow, 3).Style = "Normal"

    $ws.Cells.Item($newRow, 4).Value = "MarketMaking"
    $ws.Cells.Item($newRow, 5).Value = "UP"
    $ws.Cells.Item($newRow, 6).Value = 0.591716
    $ws.Cells.Item($newRow, 7).Value = 0.6
    $ws.Cells.Item($newRow, 8).Value = "CLOSED"
    $ws.Cells.Item($newRow, 9).Value = 1.4
    $ws.Cells.Item($newRow, 10).Value = 0.01
    $ws.Cells.Item($newRow, 11).Value = 99.44
    $ws.Cells.Item($newRow, 12).Value = 0
    $ws.Cells.Item($newRow, 13).Value = 0
    $ws.Cells.Item($newRow, 14).Value = 0.6
    $ws.Cells.Item($newRow, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($newRow, 16).Value = "early_exit"
    $ws.Cells.Item($newRow, 17).Value = 0.14
}
